$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (Week total update)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 309
$wsOff.Range("C2").Value = 209
$wsOff.Range("D2").Value = 163
$wsOff.Range("E2").Value = 63

# DEF sheet - row 2 (Week total update)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 365
$wsDef.Range("C2").Value = 268
$wsDef.Range("D2").Value = 97
$wsDef.Range("E2").Value = 49
$wsDef.Range("F2").Value = 6
